$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: SouthKorea / All (only D changes)
$ws.Cells.Item(2, 4).Value = 0.02225547035720965

# Row 3: China / All
$ws.Cells.Item(3, 5).Value = -0.0006447803591227291
$ws.Cells.Item(3, 6).Value = -0.001919317760321115
$ws.Cells.Item(3, 7).Value = 0.001274537401198386
$ws.Cells.Item(3, 8).Value = 0.6009407638284965
$ws.Cells.Item(3, 9).Value = 0.3990592361715036

# Row 4: Germany / All
$ws.Cells.Item(4, 5).Value = -0.01581963981459837
$ws.Cells.Item(4, 6).Value = -0.01392785055554586
$ws.Cells.Item(4, 7).Value = -0.00189178925905251
$ws.Cells.Item(4, 8).Value = 0.8804151496984927
$ws.Cells.Item(4, 9).Value = 0.1195848503015073

# Row 5: USA / All
$ws.Cells.Item(5, 4).Value = 0.0559832330827703
$ws.Cells.Item(5, 5).Value = -0.03372776272556065
$ws.Cells.Item(5, 6).Value = -0.01446683862051781
$ws.Cells.Item(5, 7).Value = -0.01926092410504284
$ws.Cells.Item(5, 8).Value = 0.4289296843740568
$ws.Cells.Item(5, 9).Value = 0.5710703156259431

# Row 6: was Spain/All, now becomes USA / NYC (new row inserted before old Spain/Italy)
$ws.Cells.Item(6, 1).Value = "USA"
$ws.Cells.Item(6, 2).Value = "NYC"
$ws.Cells.Item(6, 3).Value = 43943
$ws.Cells.Item(6, 3).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(6, 4).Value = 0.07258978752642781
$ws.Cells.Item(6, 5).Value = -0.05033431716921816
$ws.Cells.Item(6, 6).Value = -0.01322810364166369
$ws.Cells.Item(6, 7).Value = -0.03710621352755447
$ws.Cells.Item(6, 8).Value = 0.2628048692344893
$ws.Cells.Item(6, 9).Value = 0.7371951307655107

# Row 7: was Italy/All, now becomes Spain / All
$ws.Cells.Item(7, 1).Value = "Spain"
$ws.Cells.Item(7, 2).Value = "All"
$ws.Cells.Item(7, 3).Value = 43943
$ws.Cells.Item(7, 3).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(7, 4).Value = 0.1040126934054379
$ws.Cells.Item(7, 5).Value = -0.08175722304822823
$ws.Cells.Item(7, 6).Value = -0.0565771448318796
$ws.Cells.Item(7, 7).Value = -0.02518007821634864
$ws.Cells.Item(7, 8).Value = 0.6920140230118258
$ws.Cells.Item(7, 9).Value = 0.3079859769881741

# Row 8: new row - Italy / All
$ws.Cells.Item(8, 1).Value = "Italy"
$ws.Cells.Item(8, 2).Value = "All"
$ws.Cells.Item(8, 3).Value = 43943
$ws.Cells.Item(8, 3).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(8, 4).Value = 0.1300061819994451
$ws.Cells.Item(8, 5).Value = -0.1077507116422355
$ws.Cells.Item(8, 6).Value = -0.0707498030982413
$ws.Cells.Item(8, 7).Value = -0.03700090854399416
$ws.Cells.Item(8, 8).Value = 0.6566063650062172
$ws.Cells.Item(8, 9).Value = 0.343393634993783
